# Updated code for Create Account and Sign In page
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "Create_Account" TSID to "CreateAccount"
$ws.Range("A2").Value = "CreateAccount"

# CreateAccount test is now done (Runmode flips from N to Y)
$ws.Range("C2").Value = "Y"

# These Runmode flags flip from Y back to N
$ws.Range("C11").Value = "N"
$ws.Range("C24").Value = "N"
$ws.Range("C27").Value = "N"

# Move the active selection to F10
$ws.Range("F10").Select()
